$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.529.70'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").Value = '3.615.57'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.91'
$ws.Range("E5").Value = '  -1.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '191.47'
$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").Value = '3.611.17'
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("E8").Value = '  -2.00%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.185'
$ws.Range("E10").Value = '  +3.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.666'
$ws.Range("E11").Value = '  -0.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '56.29'
$ws.Range("E12").Value = '  -4.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000315'
$ws.Range("E13").Value = '  +8.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.76'
$ws.Range("E14").Value = '  -1.60%  '

$ws.Range("D15").Value = '4.191.88'
$ws.Range("E15").Value = '  -0.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.05'
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").Value = '3.607.44'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").Value = '70.476.71'
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.74'
$ws.Range("E19").Value = '  +0.74%  '

$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("E21").Value = '  -0.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.21'
$ws.Range("E22").Value = '  -0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.94'
$ws.Range("E23").Value = '  +3.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.94'
$ws.Range("E24").Value = '  -7.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.16'
$ws.Range("E25").Value = '  +6.63%  '

$ws.Range("E26").Value = '  -1.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.00'
$ws.Range("E27").Value = '  -4.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.13'
$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.51'
$ws.Range("E30").Value = '  -1.80%  '

$ws.Range("E31").Value = '  -1.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.29'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.120'
$ws.Range("E33").Value = '  +0.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '66.55'
$ws.Range("E34").Value = '  +1.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '580.32'
$ws.Range("E35").Value = '  -9.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '39.19'
$ws.Range("E36").Value = '  +0.96%  '

$ws.Range("D37").Value = '0.0₃0822'
$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.401'
$ws.Range("E39").Value = '  -1.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.29'
$ws.Range("E40").Value = '  +20.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  +5.79%  '

$ws.Range("E42").Value = '  -2.46%  '

$ws.Range("E43").Value = '  -6.68%  '

$ws.Range("D44").Value = '3.233.29'
$ws.Range("E44").Value = '  -2.27%  '

$ws.Range("E45").Value = '  -2.13%  '

$ws.Range("E46").Value = '  -1.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.71'
$ws.Range("E47").Value = '  +6.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.38'
$ws.Range("E48").Value = '  +3.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.28'
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  +0.10%  '
